# Auto-generated Excel COM-interop script applying the crypto price/volume refresh diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.288.84'
$ws.Range('E2').Value = '  -2.73%  '
$ws.Range('D3').Value = '1.938.45'
$ws.Range('E3').Value = '  -1.33%  '
$ws.Range('D4').Value = '''1.013'
$ws.Range('E4').Value = '  +0.71%  '
$ws.Range('D5').Value = '''321.43'
$ws.Range('E5').Value = '  -2.00%  '
$ws.Range('E6').Value = '  +0.72%  '
$ws.Range('D7').Value = '''0.4767'
$ws.Range('E7').Value = '  -4.70%  '
$ws.Range('D8').Value = '''0.4061'
$ws.Range('E8').Value = '  -3.73%  '
$ws.Range('D9').Value = '''53.28'
$ws.Range('E9').Value = '  +0.75%  '
$ws.Range('D10').Value = '''0.08482'
$ws.Range('E10').Value = '  -7.50%  '
$ws.Range('D11').Value = '''1.052'
$ws.Range('E11').Value = '  -4.41%  '
$ws.Range('D12').Value = '''22.27'
$ws.Range('E12').Value = '  -3.16%  '
$ws.Range('D13').Value = '1.984.01'
$ws.Range('E13').Value = '  -1.49%  '
$ws.Range('D14').Value = '''7.531'
$ws.Range('E14').Value = '  -4.22%  '
$ws.Range('D15').Value = '''6.119'
$ws.Range('E15').Value = '  -5.00%  '
$ws.Range('E16').Value = '  +0.94%  '
$ws.Range('D17').Value = '''89.90'
$ws.Range('E17').Value = '  -1.75%  '
$ws.Range('D18').Value = '''0.00001070'
$ws.Range('E18').Value = '  -2.82%  '
$ws.Range('D19').Value = '''0.06644'
$ws.Range('E19').Value = '  -0.73%  '
$ws.Range('D20').Value = '''18.24'
$ws.Range('E20').Value = '  -5.55%  '
$ws.Range('D21').Value = '''1.013'
$ws.Range('E21').Value = '  +0.79%  '
$ws.Range('D22').Value = '''5.815'
$ws.Range('E22').Value = '  -2.46%  '
$ws.Range('D23').Value = '28.369.52'
$ws.Range('E23').Value = '  -2.58%  '
$ws.Range('D24').Value = '''11.40'
$ws.Range('E24').Value = '  -5.67%  '
$ws.Range('D25').Value = '''2.300'
$ws.Range('E25').Value = '  +0.67%  '
$ws.Range('D26').Value = '2.231.72'
$ws.Range('E26').Value = '  -0.59%  '
$ws.Range('D27').Value = '''155.62'
$ws.Range('E27').Value = '  -0.52%  '
$ws.Range('D28').Value = '''20.22'
$ws.Range('E28').Value = '  -1.90%  '
$ws.Range('D29').Value = '''2.165'
$ws.Range('E29').Value = '  -4.37%  '
$ws.Range('D30').Value = '''5.789'
$ws.Range('E30').Value = '  -7.19%  '
$ws.Range('D31').Value = '''123.76'
$ws.Range('E31').Value = '  -2.18%  '
$ws.Range('D32').Value = '''0.9790'
$ws.Range('E32').Value = '  -6.48%  '
$ws.Range('D33').Value = '''0.09614'
$ws.Range('E33').Value = '  -2.50%  '
$ws.Range('D34').Value = '''3.682'
$ws.Range('E34').Value = '  +0.04%  '
$ws.Range('D35').Value = '''1.442'
$ws.Range('E35').Value = '  -5.60%  '
$ws.Range('D36').Value = '''5.609'
$ws.Range('E36').Value = '  -3.09%  '
$ws.Range('D37').Value = '''9.134'
$ws.Range('E37').Value = '  +1.07%  '
$ws.Range('D38').Value = '''0.02320'
$ws.Range('E38').Value = '  -4.61%  '
$ws.Range('D39').Value = '''0.06169'
$ws.Range('E39').Value = '  -3.06%  '
$ws.Range('E40').Value = '  -4.67%  '
$ws.Range('D41').Value = '''0.6196'
$ws.Range('E41').Value = '  -3.91%  '
$ws.Range('D42').Value = '''11.16'
$ws.Range('E42').Value = '  -2.47%  '
$ws.Range('D43').Value = '''1.011'
$ws.Range('E43').Value = '  +0.59%  '
$ws.Range('D44').Value = '''0.1912'
$ws.Range('E44').Value = '  -3.92%  '
$ws.Range('D45').Value = '''1.328'
$ws.Range('E45').Value = '  +3.12%  '
$ws.Range('D46').Value = '''0.5925'
$ws.Range('E46').Value = '  -5.26%  '
$ws.Range('D47').Value = '''12.81'
$ws.Range('E47').Value = '  -4.39%  '
$ws.Range('D48').Value = '''2.048'
$ws.Range('E48').Value = '  -6.73%  '
$ws.Range('E49').Value = '  -2.06%  '
$ws.Range('D50').Value = '''0.06792'
$ws.Range('E50').Value = '  -2.80%  '
$ws.Range('B51').Value = 'Quant'
$ws.Range('C51').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D51').Value = '''110.10'
$ws.Range('E51').Value = '  -2.24%  '
